$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 97, pushing rows 97-117 down to 98-118
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new weekly record
$ws.Range("A97").Value = 8
$ws.Range("B97").Value = "Terminal La Palmera de La Serena"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44588
$ws.Range("E97").Value = 4
$ws.Range("F97").Value = 100112040
$ws.Range("G97").Value = "Cilantro"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 2500
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = 2750
$ws.Range("N97").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O97").Value = "Provincia del Elquí"
$ws.Range("P97").Value = 1833
$ws.Range("Q97").Value = 1.5
$ws.Range("R97").Value = "Hortaliza"
